# Apply crypto price/volume/date updates per commit: "Updated symbol list on Thu Dec 22 22:10:04 UTC 2022 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.37"
$ws.Range("G2").Value = "'22"
$ws.Range("D3").Value = "'22.04"
$ws.Range("G3").Value = "'22"
$ws.Range("G4").Value = "'22"
$ws.Range("D5").Value = "'0.05775"
$ws.Range("G5").Value = "'22"
$ws.Range("D6").Value = "'3.421"
$ws.Range("G6").Value = "'22"
$ws.Range("D7").Value = "'6.345"
$ws.Range("G7").Value = "'22"
$ws.Range("D8").Value = "'0.8195"
$ws.Range("G8").Value = "'22"
$ws.Range("D9").Value = "'1.039"
$ws.Range("E9").Value = "'8FTXTokenFTT"
$ws.Range("G9").Value = "'22"
$ws.Range("B10").Value = "'One"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01100"
$ws.Range("E10").Value = "'9OneONEBestin24h"
$ws.Range("G10").Value = "'22"
$ws.Range("B11").Value = "'WazirX"
$ws.Range("C11").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1429"
$ws.Range("E11").Value = "'10WazirXWRX"
$ws.Range("G11").Value = "'22"
$ws.Range("B12").Value = "'MandalaExchangeToken"
$ws.Range("C12").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07300"
$ws.Range("E12").Value = "'11MandalaExchangeTokenMDX"
$ws.Range("G12").Value = "'22"
$ws.Range("B13").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03123"
$ws.Range("E13").Value = "'12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G13").Value = "'22"
$ws.Range("B14").Value = "'BitrueCoin"
$ws.Range("C14").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03117"
$ws.Range("E14").Value = "'13BitrueCoinBTR"
$ws.Range("G14").Value = "'22"
$ws.Range("B15").Value = "'MCDex"
$ws.Range("C15").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'4.145"
$ws.Range("E15").Value = "'14MCDexMCB"
$ws.Range("G15").Value = "'22"
$ws.Range("B16").Value = "'BitMartToken"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D16").Value = "'0.09384"
$ws.Range("E16").Value = "'15BitMartTokenBMX"
$ws.Range("G16").Value = "'22"
$ws.Range("B17").Value = "'BitForexToken"
$ws.Range("C17").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001595"
$ws.Range("E17").Value = "'16BitForexTokenBF"
$ws.Range("G17").Value = "'22"
$ws.Range("B18").Value = "'CoinExToken"
$ws.Range("C18").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04830"
$ws.Range("E18").Value = "'17CoinExTokenCET"
$ws.Range("G18").Value = "'22"
$ws.Range("D19").Value = "'0.006310"
$ws.Range("G19").Value = "'22"
$ws.Range("D20").Value = "'0.004133"
$ws.Range("G20").Value = "'22"
$ws.Range("D21").Value = "'0.0009921"
$ws.Range("G21").Value = "'22"
$ws.Range("D22").Value = "'0.0001498"
$ws.Range("G22").Value = "'22"
$ws.Range("D23").Value = "'3.751"
$ws.Range("G23").Value = "'22"
$ws.Range("D24").Value = "'2.191"
$ws.Range("G24").Value = "'22"
$ws.Range("D25").Value = "'0.3234"
$ws.Range("G25").Value = "'22"
$ws.Range("G26").Value = "'22"
$ws.Range("D27").Value = "'0.0003995"
$ws.Range("G27").Value = "'22"
$ws.Range("G28").Value = "'22"
$ws.Range("G29").Value = "'22"
$ws.Range("G30").Value = "'22"
$ws.Range("G31").Value = "'22"
$ws.Range("G32").Value = "'22"
$ws.Range("G33").Value = "'22"
$ws.Range("G34").Value = "'22"
$ws.Range("G35").Value = "'22"
$ws.Range("G36").Value = "'22"
$ws.Range("G37").Value = "'22"
$ws.Range("G38").Value = "'22"
$ws.Range("G39").Value = "'22"
$ws.Range("D40").Value = "'0.03869"
$ws.Range("G40").Value = "'22"
$ws.Range("D41").Value = "'0.006669"
$ws.Range("G41").Value = "'22"
$ws.Range("D42").Value = "'0.1071"
$ws.Range("G42").Value = "'22"
$ws.Range("D43").Value = "'0.002897"
$ws.Range("G43").Value = "'22"
$ws.Range("D44").Value = "'0.006580"
$ws.Range("G44").Value = "'22"
$ws.Range("D45").Value = "'0.00005590"
$ws.Range("G45").Value = "'22"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("G46").Value = "'22"
$ws.Range("D47").Value = "'0.3896"
$ws.Range("G47").Value = "'22"
$ws.Range("G48").Value = "'22"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("G49").Value = "'22"
$ws.Range("G50").Value = "'22"
$ws.Range("G51").Value = "'22"
